# "excel actualizados a GKD only"
# Removes the MDG-*/SOM-* instance rows (29-33) from the results sheet,
# leaving only the GKD-* instances. The dependent AVERAGE() formulas in
# rows 3-6 and the shared-strings table recalculate/compact automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the instance name (col A) and all metric columns (B:AC) for the
# five rows that correspond to the MDG-a, MDG-b, SOM-a (x2) and SOM-b
# instances (rows 29-33), while keeping the pre-existing cell formatting.
$ws.Range("A29:AC33").ClearContents()

# Force recalculation so the AVERAGE() formulas in rows 3-6 reflect the
# remaining GKD-only data.
$excel.CalculateFullRebuild()

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("I14").Select()
